$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$defaultStyle = $ws.Range("B2").Style

$ws.Range("D2").Value = "'29.055.42"
$ws.Range("D2").Style = $defaultStyle
$ws.Range("E2").Value = "  +0.09%  "

$ws.Range("D3").Value = "'1.829.35"
$ws.Range("D3").Style = $defaultStyle
$ws.Range("E3").Value = "  -0.01%  "

$ws.Range("D4").Value = "'0.9985"
$ws.Range("D4").Style = $defaultStyle
$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").Value = "'244.63"
$ws.Range("D5").Style = $defaultStyle
$ws.Range("E5").Value = "  +1.34%  "

$ws.Range("D6").Value = "'0.6330"
$ws.Range("D6").Style = $defaultStyle
$ws.Range("E6").Value = "  +0.88%  "

$ws.Range("D7").Value = "'1.000"
$ws.Range("D7").Style = $defaultStyle
$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("D8").Value = "'0.07533"
$ws.Range("D8").Style = $defaultStyle
$ws.Range("E8").Value = "  -0.89%  "

$ws.Range("D9").Value = "'0.2947"
$ws.Range("D9").Style = $defaultStyle
$ws.Range("E9").Value = "  +1.08%  "

$ws.Range("D10").Value = "'23.14"
$ws.Range("D10").Style = $defaultStyle
$ws.Range("E10").Value = "  +1.54%  "

$ws.Range("D11").Value = "'0.07708"
$ws.Range("D11").Style = $defaultStyle
$ws.Range("E11").Value = "  +0.83%  "

$ws.Range("D12").Value = "'1.831.29"
$ws.Range("D12").Style = $defaultStyle
$ws.Range("E12").Value = "  -0.22%  "

$ws.Range("D13").Value = "'5.002"
$ws.Range("D13").Style = $defaultStyle
$ws.Range("E13").Value = "  +0.91%  "

$ws.Range("D14").Value = "'0.6701"
$ws.Range("D14").Style = $defaultStyle
$ws.Range("E14").Value = "  +0.77%  "

$ws.Range("D15").Value = "'83.24"
$ws.Range("D15").Style = $defaultStyle
$ws.Range("E15").Value = "  +1.07%  "

$ws.Range("D16").Value = "'0.000009609"
$ws.Range("D16").Style = $defaultStyle
$ws.Range("E16").Value = "  +1.89%  "

$ws.Range("D17").Value = "'6.065"
$ws.Range("D17").Style = $defaultStyle
$ws.Range("E17").Value = "  +1.43%  "

$ws.Range("D18").Value = "'29.101.00"
$ws.Range("D18").Style = $defaultStyle
$ws.Range("E18").Value = "  +0.51%  "

$ws.Range("D19").Value = "'12.61"
$ws.Range("D19").Style = $defaultStyle
$ws.Range("E19").Value = "  +2.36%  "

$ws.Range("D20").Value = "'226.64"
$ws.Range("D20").Style = $defaultStyle
$ws.Range("E20").Value = "  +0.82%  "

$ws.Range("D21").Value = "'0.9999"
$ws.Range("D21").Style = $defaultStyle
$ws.Range("E21").Value = "  +0.04%  "

$ws.Range("D22").Value = "'7.157"
$ws.Range("D22").Style = $defaultStyle
$ws.Range("E22").Value = "  -1.06%  "

$ws.Range("D23").Value = "'0.9995"
$ws.Range("D23").Style = $defaultStyle
$ws.Range("E23").Value = "  -0.12%  "

$ws.Range("D24").Value = "'160.63"
$ws.Range("D24").Style = $defaultStyle
$ws.Range("E24").Value = "  -0.25%  "

$ws.Range("D25").Value = "'0.1425"
$ws.Range("D25").Style = $defaultStyle
$ws.Range("E25").Value = "  +4.31%  "

$ws.Range("E26").Value = "  +1.44%  "

$ws.Range("D27").Value = "'17.97"
$ws.Range("D27").Style = $defaultStyle
$ws.Range("E27").Value = "  +0.79%  "

$ws.Range("D28").Value = "'1.507"
$ws.Range("D28").Style = $defaultStyle
$ws.Range("E28").Value = "  +0.72%  "

$ws.Range("D29").Value = "'4.150"
$ws.Range("D29").Style = $defaultStyle
$ws.Range("E29").Value = "  +2.43%  "

$ws.Range("D30").Value = "'4.067"
$ws.Range("D30").Style = $defaultStyle
$ws.Range("E30").Value = "  +0.93%  "

$ws.Range("D31").Value = "'0.05490"
$ws.Range("D31").Style = $defaultStyle
$ws.Range("E31").Value = "  +5.52%  "

$ws.Range("D32").Value = "'1.203"
$ws.Range("D32").Style = $defaultStyle
$ws.Range("E32").Value = "  +0.05%  "

$ws.Range("D33").Value = "'1.861"
$ws.Range("D33").Style = $defaultStyle
$ws.Range("E33").Value = "  +0.62%  "

$ws.Range("D34").Value = "'0.7462"
$ws.Range("D34").Style = $defaultStyle
$ws.Range("E34").Value = "  +2.24%  "

$ws.Range("D35").Value = "'1.140"
$ws.Range("D35").Style = $defaultStyle
$ws.Range("E35").Value = "  -1.11%  "

$ws.Range("D36").Value = "'2.652"
$ws.Range("D36").Style = $defaultStyle
$ws.Range("E36").Value = "  +1.56%  "

$ws.Range("D37").Value = "'1.244.77"
$ws.Range("D37").Style = $defaultStyle
$ws.Range("E37").Value = "  -2.19%  "

$ws.Range("D38").Value = "'2.753"
$ws.Range("D38").Style = $defaultStyle
$ws.Range("E38").Value = "  -0.14%  "

$ws.Range("D39").Value = "'0.01785"
$ws.Range("D39").Style = $defaultStyle
$ws.Range("E39").Value = "  +0.00%  "

$ws.Range("D40").Value = "'6.637"
$ws.Range("D40").Style = $defaultStyle
$ws.Range("E40").Value = "  +2.47%  "

$ws.Range("D41").Value = "'0.9039"
$ws.Range("D41").Style = $defaultStyle
$ws.Range("E41").Value = "  +1.27%  "

$ws.Range("D42").Value = "'1.000"
$ws.Range("D42").Style = $defaultStyle
$ws.Range("E42").Value = "  +0.01%  "

$ws.Range("D43").Value = "'101.47"
$ws.Range("D43").Style = $defaultStyle
$ws.Range("E43").Value = "  -0.08%  "

$ws.Range("D44").Value = "'1.974.51"
$ws.Range("D44").Style = $defaultStyle
$ws.Range("E44").Value = "  -0.03%  "

$ws.Range("D45").Value = "'65.25"
$ws.Range("D45").Style = $defaultStyle
$ws.Range("E45").Value = "  +2.40%  "

$ws.Range("E46").Value = "  +1.05%  "

$ws.Range("D47").Value = "'0.5096"
$ws.Range("D47").Style = $defaultStyle
$ws.Range("E47").Value = "  -0.21%  "

$ws.Range("D48").Value = "'0.4064"
$ws.Range("D48").Style = $defaultStyle
$ws.Range("E48").Value = "  +2.06%  "

$ws.Range("D49").Value = "'9.015"
$ws.Range("D49").Style = $defaultStyle
$ws.Range("E49").Value = "  +1.88%  "

$ws.Range("D50").Value = "'1.660"
$ws.Range("D50").Style = $defaultStyle
$ws.Range("E50").Value = "  +0.96%  "

$ws.Range("D51").Value = "'6.786"
$ws.Range("D51").Style = $defaultStyle
$ws.Range("E51").Value = "  +1.35%  "
